# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price updates to the Alpha_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 591
$ws.Range("I4").Value = 591
$ws.Range("K4").Value = 591
$ws.Range("M4").Value = -477

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 8416.5
$ws.Range("J32").Value = 8199.9
$ws.Range("L32").Value = 8199.9
$ws.Range("N32").Value = -8851.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2546.889
$ws.Range("I98").Value = 2309.75
$ws.Range("J98").Value = 4444
$ws.Range("K98").Value = 2309.75
$ws.Range("L98").Value = 4444
$ws.Range("M98").Value = -811.75
$ws.Range("N98").Value = -7440

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 4999
$ws.Range("I111").Value = 4999
$ws.Range("K111").Value = 14997
$ws.Range("M111").Value = -11930

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2520.9412
$ws.Range("I113").Value = 2453.1428
$ws.Range("K113").Value = 2453.1428
$ws.Range("M113").Value = 800.8571999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2546.889
$ws.Range("I122").Value = 2309.75
$ws.Range("J122").Value = 4444
$ws.Range("K122").Value = 6929.25
$ws.Range("L122").Value = 13332
$ws.Range("M122").Value = -4479.25
$ws.Range("N122").Value = -18232

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1328.3429
$ws.Range("J138").Value = 2723.125
$ws.Range("L138").Value = 8169.375
$ws.Range("N138").Value = -18449.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2215.5
$ws.Range("I32").Value = 2230.2593
$ws.Range("K32").Value = 2230.2593
$ws.Range("M32").Value = -1943.2593

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2132.5
$ws.Range("I122").Value = 1866.2858
$ws.Range("K122").Value = 5598.857400000001
$ws.Range("M122").Value = -3148.857400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3274.2778
$ws.Range("I86").Value = 3438.6667
$ws.Range("K86").Value = 3438.6667
$ws.Range("M86").Value = -2315.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3274.2778
$ws.Range("I89").Value = 3438.6667
$ws.Range("K89").Value = 17193.3335
$ws.Range("M89").Value = -11577.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2167.7778
$ws.Range("I94").Value = 2252
$ws.Range("K94").Value = 2252
$ws.Range("M94").Value = -1801

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 579.5
$ws.Range("I22").Value = 616.5714
$ws.Range("K22").Value = 616.5714
$ws.Range("M22").Value = -266.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 8557.571
$ws.Range("I23").Value = 5817.1665
$ws.Range("J23").Value = 25000
$ws.Range("K23").Value = 5817.1665
$ws.Range("L23").Value = 25000
$ws.Range("M23").Value = -5577.1665
$ws.Range("N23").Value = -25480

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H27").Value = 8557.571
$ws.Range("I27").Value = 5817.1665
$ws.Range("J27").Value = 25000
$ws.Range("K27").Value = 5817.1665
$ws.Range("L27").Value = 25000
$ws.Range("M27").Value = -5625.1665
$ws.Range("N27").Value = -25384

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2209.111
$ws.Range("I122").Value = 1849.75
$ws.Range("K122").Value = 5549.25
$ws.Range("M122").Value = -3099.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 478.91666
$ws.Range("J7").Value = 423.14285
$ws.Range("L7").Value = 1269.42855
$ws.Range("N7").Value = -1493.42855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 298.85715
$ws.Range("I17").Value = 132.66667
$ws.Range("J17").Value = 423.5
$ws.Range("K17").Value = 398.00001
$ws.Range("L17").Value = 1270.5
$ws.Range("M17").Value = -229.00001
$ws.Range("N17").Value = -1608.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 765.8570999999999
$ws.Range("I33").Value = 455
$ws.Range("J33").Value = 890.2
$ws.Range("K33").Value = 2730
$ws.Range("L33").Value = 5341.200000000001
$ws.Range("M33").Value = -2447
$ws.Range("N33").Value = -5907.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 607553.7
$ws.Range("J131").Value = 828055.2
$ws.Range("L131").Value = 2484165.6
$ws.Range("N131").Value = -2494245.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8468.105
$ws.Range("I70").Value = 8024.625
$ws.Range("K70").Value = 8024.625
$ws.Range("M70").Value = -7754.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 8468.105
$ws.Range("I73").Value = 8024.625
$ws.Range("K73").Value = 8024.625
$ws.Range("M73").Value = -7088.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7086.9165
$ws.Range("I80").Value = 4719.2856
$ws.Range("K80").Value = 4719.2856
$ws.Range("M80").Value = -3721.2856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 7086.9165
$ws.Range("I83").Value = 4719.2856
$ws.Range("K83").Value = 23596.428
$ws.Range("M83").Value = -18604.428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2718.6
$ws.Range("I102").Value = 4296.5
$ws.Range("J102").Value = 1666.6666
$ws.Range("K102").Value = 4296.5
$ws.Range("L102").Value = 1666.6666
$ws.Range("M102").Value = -2674.5
$ws.Range("N102").Value = -4910.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3055.24
$ws.Range("I122").Value = 3371.5
$ws.Range("K122").Value = 10114.5
$ws.Range("M122").Value = -7664.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1396.875
$ws.Range("I22").Value = 850
$ws.Range("J22").Value = 2100
$ws.Range("K22").Value = 850
$ws.Range("L22").Value = 2100
$ws.Range("M22").Value = -555
$ws.Range("N22").Value = -2690

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1396.875
$ws.Range("I27").Value = 850
$ws.Range("J27").Value = 2100
$ws.Range("K27").Value = 850
$ws.Range("L27").Value = 2100
$ws.Range("M27").Value = -743
$ws.Range("N27").Value = -2314

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 2490.5557
$ws.Range("I31").Value = 471.45456
$ws.Range("J31").Value = 5663.4287
$ws.Range("K31").Value = 471.45456
$ws.Range("L31").Value = 5663.4287
$ws.Range("M31").Value = -223.45456
$ws.Range("N31").Value = -6159.4287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2918.75
$ws.Range("I40").Value = 2333.077
$ws.Range("J40").Value = 5456.6665
$ws.Range("K40").Value = 2333.077
$ws.Range("L40").Value = 5456.6665
$ws.Range("M40").Value = -2197.077
$ws.Range("N40").Value = -5728.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3595.7
$ws.Range("I122").Value = 3330.6667
$ws.Range("K122").Value = 9992.000100000001
$ws.Range("M122").Value = -7542.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 12398.8
$ws.Range("I18").Value = 12398.8
$ws.Range("K18").Value = 12398.8
$ws.Range("M18").Value = -12225.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1196.091
$ws.Range("I107").Value = 807.125
$ws.Range("J107").Value = 2233.3333
$ws.Range("K107").Value = 2421.375
$ws.Range("L107").Value = 6699.999899999999
$ws.Range("M107").Value = -501.375
$ws.Range("N107").Value = -10539.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2551.4443
$ws.Range("I126").Value = 1821
$ws.Range("J126").Value = 3135.8
$ws.Range("K126").Value = 5463
$ws.Range("L126").Value = 9407.400000000001
$ws.Range("M126").Value = -2993
$ws.Range("N126").Value = -14347.4
